$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(56).Insert()
Write-Host ("Q57 before copy: ")
Write-Host $ws.Range("Q57").Value2
$ws.Range("A57:Q57").Copy()
$ws.Range("A56:Q56").PasteSpecial(-4122)
Write-Host "Done"
